$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.260118246078491
$ws.Range("B1").Value = 3.884607791900635
$ws.Range("C1").Value = 3.630458354949951
$ws.Range("D1").Value = 3.429226636886597
$ws.Range("E1").Value = 1.075165867805481
